$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows 219 and 220 (copy style for columns A and E from an existing data row)
$ws.Range("A218").Copy($ws.Range("A219"))
$ws.Range("E218").Copy($ws.Range("E219"))
$ws.Range("A218").Copy($ws.Range("A220"))
$ws.Range("E218").Copy($ws.Range("E220"))

# Row 32
$ws.Cells.Item(32, 2).Value = 5400043
$ws.Cells.Item(32, 6).Value = "Sydney FC"
$ws.Cells.Item(32, 7).Value = "Central Coast Mariners"
$ws.Cells.Item(32, 8).Value = 3
$ws.Cells.Item(32, 9).Value = 2
$ws.Cells.Item(32, 11).Value = 2.15
$ws.Cells.Item(32, 12).Value = 3.4
$ws.Cells.Item(32, 13).Value = 3.4
$ws.Cells.Item(32, 14).Value = 2.2
$ws.Cells.Item(32, 15).Value = 3.75
$ws.Cells.Item(32, 16).Value = 3
$ws.Cells.Item(32, 17).Value = -0.25
$ws.Cells.Item(32, 18).Value = 1.975
$ws.Cells.Item(32, 19).Value = 1.875
$ws.Cells.Item(32, 20).Value = 3
$ws.Cells.Item(32, 21).Value = 1.95
$ws.Cells.Item(32, 22).Value = 1.9
$ws.Cells.Item(32, 23).Value = 1.2
$ws.Cells.Item(32, 26).Value = 0.9750000000000001
$ws.Cells.Item(32, 28).Value = 0.95

# Row 33
$ws.Cells.Item(33, 2).Value = 5400042
$ws.Cells.Item(33, 6).Value = "Melbourne City"
$ws.Cells.Item(33, 7).Value = "Macarthur FC"
$ws.Cells.Item(33, 8).Value = 6
$ws.Cells.Item(33, 9).Value = 1
$ws.Cells.Item(33, 11).Value = 1.533
$ws.Cells.Item(33, 12).Value = 4.2
$ws.Cells.Item(33, 13).Value = 6
$ws.Cells.Item(33, 14).Value = 1.333
$ws.Cells.Item(33, 15).Value = 5
$ws.Cells.Item(33, 16).Value = 9.5
$ws.Cells.Item(33, 17).Value = -1.5
$ws.Cells.Item(33, 18).Value = 1.85
$ws.Cells.Item(33, 19).Value = 2
$ws.Cells.Item(33, 20).Value = 3.25
$ws.Cells.Item(33, 21).Value = 1.875
$ws.Cells.Item(33, 22).Value = 1.975
$ws.Cells.Item(33, 23).Value = 0.333
$ws.Cells.Item(33, 26).Value = 0.8500000000000001
$ws.Cells.Item(33, 28).Value = 0.875

# Row 97
$ws.Cells.Item(97, 2).Value = 5400063
$ws.Cells.Item(97, 6).Value = "Melbourne City"
$ws.Cells.Item(97, 7).Value = "Western Sydney Wanderers"
$ws.Cells.Item(97, 8).Value = 3
$ws.Cells.Item(97, 9).Value = 2
$ws.Cells.Item(97, 10).Value = "H"
$ws.Cells.Item(97, 11).Value = 1.75
$ws.Cells.Item(97, 12).Value = 3.8
$ws.Cells.Item(97, 13).Value = 4
$ws.Cells.Item(97, 14).Value = 2
$ws.Cells.Item(97, 16).Value = 3.4
$ws.Cells.Item(97, 17).Value = -0.5
$ws.Cells.Item(97, 18).Value = 2.025
$ws.Cells.Item(97, 19).Value = 1.825
$ws.Cells.Item(97, 20).Value = 3
$ws.Cells.Item(97, 21).Value = 1.85
$ws.Cells.Item(97, 22).Value = 2
$ws.Cells.Item(97, 23).Value = 1
$ws.Cells.Item(97, 25).Value = -1
$ws.Cells.Item(97, 26).Value = 1.025
$ws.Cells.Item(97, 27).Value = -1
$ws.Cells.Item(97, 28).Value = 0.8500000000000001

# Row 98
$ws.Cells.Item(98, 2).Value = 5404732
$ws.Cells.Item(98, 6).Value = "Adelaide United"
$ws.Cells.Item(98, 7).Value = "Central Coast Mariners"
$ws.Cells.Item(98, 8).Value = 1
$ws.Cells.Item(98, 9).Value = 4
$ws.Cells.Item(98, 10).Value = "A"
$ws.Cells.Item(98, 11).Value = 2.3
$ws.Cells.Item(98, 12).Value = 3.75
$ws.Cells.Item(98, 13).Value = 2.75
$ws.Cells.Item(98, 14).Value = 2.625
$ws.Cells.Item(98, 16).Value = 2.4
$ws.Cells.Item(98, 17).Value = 0
$ws.Cells.Item(98, 18).Value = 2.05
$ws.Cells.Item(98, 19).Value = 1.8
$ws.Cells.Item(98, 20).Value = 3.75
$ws.Cells.Item(98, 21).Value = 2
$ws.Cells.Item(98, 22).Value = 1.85
$ws.Cells.Item(98, 23).Value = -1
$ws.Cells.Item(98, 25).Value = 1.4
$ws.Cells.Item(98, 26).Value = -1
$ws.Cells.Item(98, 27).Value = 0.8
$ws.Cells.Item(98, 28).Value = 1

# Row 99
$ws.Cells.Item(99, 2).Value = 5404735
$ws.Cells.Item(99, 6).Value = "Macarthur FC"
$ws.Cells.Item(99, 7).Value = "Wellington Phoenix"
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 1
$ws.Cells.Item(99, 10).Value = "A"
$ws.Cells.Item(99, 11).Value = 3.6
$ws.Cells.Item(99, 12).Value = 3.75
$ws.Cells.Item(99, 13).Value = 1.909
$ws.Cells.Item(99, 14).Value = 4
$ws.Cells.Item(99, 15).Value = 4
$ws.Cells.Item(99, 16).Value = 1.833
$ws.Cells.Item(99, 17).Value = 0.5
$ws.Cells.Item(99, 20).Value = 3.25
$ws.Cells.Item(99, 21).Value = 1.925
$ws.Cells.Item(99, 22).Value = 1.925
$ws.Cells.Item(99, 23).Value = -1
$ws.Cells.Item(99, 25).Value = 0.833
$ws.Cells.Item(99, 26).Value = -1
$ws.Cells.Item(99, 27).Value = 0.825
$ws.Cells.Item(99, 29).Value = 0.925

# Row 100
$ws.Cells.Item(100, 2).Value = 5400064
$ws.Cells.Item(100, 6).Value = "Sydney FC"
$ws.Cells.Item(100, 7).Value = "Newcastle Jets"
$ws.Cells.Item(100, 8).Value = 2
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = "H"
$ws.Cells.Item(100, 11).Value = 1.65
$ws.Cells.Item(100, 12).Value = 4
$ws.Cells.Item(100, 13).Value = 4.5
$ws.Cells.Item(100, 14).Value = 1.533
$ws.Cells.Item(100, 15).Value = 4.75
$ws.Cells.Item(100, 16).Value = 5.5
$ws.Cells.Item(100, 17).Value = -1.25
$ws.Cells.Item(100, 20).Value = 3.5
$ws.Cells.Item(100, 21).Value = 2
$ws.Cells.Item(100, 22).Value = 1.85
$ws.Cells.Item(100, 23).Value = 0.5329999999999999
$ws.Cells.Item(100, 25).Value = -1
$ws.Cells.Item(100, 26).Value = 1.025
$ws.Cells.Item(100, 27).Value = -1
$ws.Cells.Item(100, 29).Value = 0.8500000000000001

# Row 180
$ws.Cells.Item(180, 2).Value = 7646749
$ws.Cells.Item(180, 6).Value = "Brisbane Roar"
$ws.Cells.Item(180, 7).Value = "Newcastle Jets"
$ws.Cells.Item(180, 9).Value = 2
$ws.Cells.Item(180, 10).Value = "H"
$ws.Cells.Item(180, 11).Value = 1.909
$ws.Cells.Item(180, 12).Value = 4
$ws.Cells.Item(180, 13).Value = 3.4
$ws.Cells.Item(180, 14).Value = 2.4
$ws.Cells.Item(180, 15).Value = 4
$ws.Cells.Item(180, 16).Value = 2.6
$ws.Cells.Item(180, 17).Value = 0
$ws.Cells.Item(180, 18).Value = 1.83
$ws.Cells.Item(180, 19).Value = 2.07
$ws.Cells.Item(180, 20).Value = 3.25
$ws.Cells.Item(180, 21).Value = 1.9
$ws.Cells.Item(180, 22).Value = 1.95
$ws.Cells.Item(180, 23).Value = 1.4
$ws.Cells.Item(180, 25).Value = -1
$ws.Cells.Item(180, 26).Value = 0.8300000000000001
$ws.Cells.Item(180, 27).Value = -1
$ws.Cells.Item(180, 28).Value = 0.8999999999999999

# Row 181
$ws.Cells.Item(181, 2).Value = 7646750
$ws.Cells.Item(181, 6).Value = "Perth Glory"
$ws.Cells.Item(181, 7).Value = "Wellington Phoenix"
$ws.Cells.Item(181, 9).Value = 4
$ws.Cells.Item(181, 10).Value = "A"
$ws.Cells.Item(181, 11).Value = 2.45
$ws.Cells.Item(181, 12).Value = 3.75
$ws.Cells.Item(181, 13).Value = 2.55
$ws.Cells.Item(181, 14).Value = 3.1
$ws.Cells.Item(181, 15).Value = 3.8
$ws.Cells.Item(181, 16).Value = 2.05
$ws.Cells.Item(181, 17).Value = 0.25
$ws.Cells.Item(181, 18).Value = 2
$ws.Cells.Item(181, 19).Value = 1.85
$ws.Cells.Item(181, 20).Value = 3
$ws.Cells.Item(181, 21).Value = 1.925
$ws.Cells.Item(181, 22).Value = 1.925
$ws.Cells.Item(181, 23).Value = -1
$ws.Cells.Item(181, 25).Value = 1.05
$ws.Cells.Item(181, 26).Value = -1
$ws.Cells.Item(181, 27).Value = 0.8500000000000001
$ws.Cells.Item(181, 28).Value = 0.925

# Row 214
$ws.Cells.Item(214, 2).Value = 7127375
$ws.Cells.Item(214, 5).Value = 45345.23958333334
$ws.Cells.Item(214, 6).Value = "Brisbane Roar"
$ws.Cells.Item(214, 7).Value = "Western United FC"
$ws.Cells.Item(214, 8).Value = 2
$ws.Cells.Item(214, 9).Value = 2
$ws.Cells.Item(214, 10).Value = "D"
$ws.Cells.Item(214, 11).Value = 2
$ws.Cells.Item(214, 12).Value = 3.75
$ws.Cells.Item(214, 13).Value = 3.2
$ws.Cells.Item(214, 14).Value = 1.727
$ws.Cells.Item(214, 15).Value = 4.333
$ws.Cells.Item(214, 16).Value = 4
$ws.Cells.Item(214, 17).Value = -0.75
$ws.Cells.Item(214, 18).Value = 1.925
$ws.Cells.Item(214, 19).Value = 1.925
$ws.Cells.Item(214, 21).Value = 1.9
$ws.Cells.Item(214, 22).Value = 1.95
$ws.Cells.Item(214, 23).Value = -1
$ws.Cells.Item(214, 24).Value = 3.333
$ws.Cells.Item(214, 25).Value = -1
$ws.Cells.Item(214, 26).Value = -1
$ws.Cells.Item(214, 27).Value = 0.925
$ws.Cells.Item(214, 28).Value = 0.8999999999999999
$ws.Cells.Item(214, 29).Value = -1

# Row 215
$ws.Cells.Item(215, 2).Value = 7875268
$ws.Cells.Item(215, 5).Value = 45345.38194444445
$ws.Cells.Item(215, 6).Value = "Perth Glory"
$ws.Cells.Item(215, 7).Value = "Sydney FC"
$ws.Cells.Item(215, 8).Value = 1
$ws.Cells.Item(215, 9).Value = 1
$ws.Cells.Item(215, 10).Value = "D"
$ws.Cells.Item(215, 11).Value = 1.75
$ws.Cells.Item(215, 12).Value = 4
$ws.Cells.Item(215, 13).Value = 3.4
$ws.Cells.Item(215, 14).Value = 1.75
$ws.Cells.Item(215, 16).Value = 3.4
$ws.Cells.Item(215, 17).Value = -0.5
$ws.Cells.Item(215, 18).Value = 1.825
$ws.Cells.Item(215, 19).Value = 2.025
$ws.Cells.Item(215, 21).Value = 1.875
$ws.Cells.Item(215, 22).Value = 1.975
$ws.Cells.Item(215, 23).Value = -1
$ws.Cells.Item(215, 24).Value = 3
$ws.Cells.Item(215, 25).Value = -1
$ws.Cells.Item(215, 26).Value = -1
$ws.Cells.Item(215, 27).Value = 1.025
$ws.Cells.Item(215, 28).Value = -1
$ws.Cells.Item(215, 29).Value = 0.9750000000000001

# Row 216
$ws.Cells.Item(216, 2).Value = 7126789
$ws.Cells.Item(216, 5).Value = 45346.14583333334
$ws.Cells.Item(216, 6).Value = "Sydney FC"
$ws.Cells.Item(216, 7).Value = "Melbourne City"
$ws.Cells.Item(216, 11).Value = 1.833
$ws.Cells.Item(216, 12).Value = 4.2
$ws.Cells.Item(216, 13).Value = 3.75
$ws.Cells.Item(216, 14).Value = 1.833
$ws.Cells.Item(216, 15).Value = 4.2
$ws.Cells.Item(216, 16).Value = 3.75
$ws.Cells.Item(216, 17).Value = -0.5
$ws.Cells.Item(216, 18).Value = 1.86
$ws.Cells.Item(216, 19).Value = 2.04
$ws.Cells.Item(216, 20).Value = 3.25
$ws.Cells.Item(216, 21).Value = 1.825
$ws.Cells.Item(216, 22).Value = 2.025

# Row 217
$ws.Cells.Item(217, 2).Value = 7127377
$ws.Cells.Item(217, 5).Value = 45346.23958333334
$ws.Cells.Item(217, 6).Value = "Adelaide United"
$ws.Cells.Item(217, 7).Value = "Western Sydney Wanderers"
$ws.Cells.Item(217, 11).Value = 2.4
$ws.Cells.Item(217, 12).Value = 4
$ws.Cells.Item(217, 13).Value = 2.6
$ws.Cells.Item(217, 14).Value = 2.4
$ws.Cells.Item(217, 16).Value = 2.6
$ws.Cells.Item(217, 17).Value = 0
$ws.Cells.Item(217, 18).Value = 1.84
$ws.Cells.Item(217, 19).Value = 2.06
$ws.Cells.Item(217, 21).Value = 1.85
$ws.Cells.Item(217, 22).Value = 2

# Row 218
$ws.Cells.Item(218, 2).Value = 7127378
$ws.Cells.Item(218, 5).Value = 45346.32291666666
$ws.Cells.Item(218, 6).Value = "Perth Glory"
$ws.Cells.Item(218, 7).Value = "Wellington Phoenix"
$ws.Cells.Item(218, 11).Value = 2.375
$ws.Cells.Item(218, 12).Value = 3.5
$ws.Cells.Item(218, 13).Value = 2.875
$ws.Cells.Item(218, 14).Value = 2.375
$ws.Cells.Item(218, 15).Value = 3.5
$ws.Cells.Item(218, 16).Value = 2.875
$ws.Cells.Item(218, 17).Value = -0.25
$ws.Cells.Item(218, 18).Value = 2.07
$ws.Cells.Item(218, 19).Value = 1.83
$ws.Cells.Item(218, 20).Value = 3
$ws.Cells.Item(218, 21).Value = 1.975
$ws.Cells.Item(218, 22).Value = 1.875

# Row 219
$ws.Cells.Item(219, 1).Value = 217
$ws.Cells.Item(219, 2).Value = 7127376
$ws.Cells.Item(219, 3).Value = "Australia ALeague"
$ws.Cells.Item(219, 4).Value = "Australia ALeague"
$ws.Cells.Item(219, 5).Value = 45347.125
$ws.Cells.Item(219, 6).Value = "Newcastle Jets"
$ws.Cells.Item(219, 7).Value = "Macarthur FC"
$ws.Cells.Item(219, 11).Value = 1.95
$ws.Cells.Item(219, 12).Value = 4
$ws.Cells.Item(219, 13).Value = 3.4
$ws.Cells.Item(219, 14).Value = 1.95
$ws.Cells.Item(219, 15).Value = 4
$ws.Cells.Item(219, 16).Value = 3.4
$ws.Cells.Item(219, 17).Value = -0.5
$ws.Cells.Item(219, 18).Value = 2
$ws.Cells.Item(219, 19).Value = 1.9
$ws.Cells.Item(219, 20).Value = 3.25
$ws.Cells.Item(219, 21).Value = 1.875
$ws.Cells.Item(219, 22).Value = 1.975
$ws.Cells.Item(219, 23).Value = 0
$ws.Cells.Item(219, 24).Value = 0
$ws.Cells.Item(219, 25).Value = 0
$ws.Cells.Item(219, 26).Value = 0
$ws.Cells.Item(219, 27).Value = 0

# Row 220
$ws.Cells.Item(220, 1).Value = 218
$ws.Cells.Item(220, 2).Value = 7127379
$ws.Cells.Item(220, 3).Value = "Australia ALeague"
$ws.Cells.Item(220, 4).Value = "Australia ALeague"
$ws.Cells.Item(220, 5).Value = 45347.125
$ws.Cells.Item(220, 6).Value = "Melbourne Victory"
$ws.Cells.Item(220, 7).Value = "Central Coast Mariners"
$ws.Cells.Item(220, 11).Value = 1.95
$ws.Cells.Item(220, 12).Value = 3.6
$ws.Cells.Item(220, 13).Value = 3.8
$ws.Cells.Item(220, 14).Value = 1.95
$ws.Cells.Item(220, 15).Value = 3.6
$ws.Cells.Item(220, 16).Value = 3.8
$ws.Cells.Item(220, 17).Value = -0.5
$ws.Cells.Item(220, 18).Value = 1.98
$ws.Cells.Item(220, 19).Value = 1.92
$ws.Cells.Item(220, 20).Value = 2.75
$ws.Cells.Item(220, 21).Value = 1.95
$ws.Cells.Item(220, 22).Value = 1.9
$ws.Cells.Item(220, 23).Value = 0
$ws.Cells.Item(220, 24).Value = 0
$ws.Cells.Item(220, 25).Value = 0
$ws.Cells.Item(220, 26).Value = 0
$ws.Cells.Item(220, 27).Value = 0
